# CubeA-HW25: two new simulation methods ("Holden", "Rizzie Spiral") were added
# to the sweep, and "Thomas Hex" was renamed to "Matthies Hex". The whole
# simulation was rerun, so every method's row of results moved down by two
# rows (to make room for the two new methods, inserted right after "Spiral5")
# and got fresh numeric results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Snapshot the existing data rows (3..29) BEFORE we overwrite anything.
#    Row index in the sheet = method index (A column) + 2 (header rows 1-2).
# ---------------------------------------------------------------------------
$lastOldRow = 29
$firstDataRow = 4   # method index 2 .. 27 today (method index 0 is row2, 1 is row3)

$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastOldRow; $r++) {
    $vals = @()
    for ($c = 3; $c -le 23; $c++) {
        $vals += ,$ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $vals
}

# ---------------------------------------------------------------------------
# 2) Full, ordered list of method names after the edit (this is column B,
#    shared-string labels, in row order starting at row 2 / method index 0).
# ---------------------------------------------------------------------------
$methods = @(
    "HKL",
    "Spiral5",
    "Holden",
    "Rizzie Spiral",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Matthies Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

# ---------------------------------------------------------------------------
# 3) Move the old data rows (4..29) down by two (-> 6..31), restoring their
#    original values unchanged.
# ---------------------------------------------------------------------------
for ($r = $lastOldRow; $r -ge $firstDataRow; $r--) {
    $newRow = $r + 2
    $vals = $snapshot[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($newRow, 3 + $i).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------------
# 4) Brand-new simulation results for the two newly added methods, which now
#    occupy rows 4 and 5 (method indices 2 and 3).
# ---------------------------------------------------------------------------
$newRow4 = @(1.132255175101699,0.953527778241232,0.9235374354316438,1.040972862200088,0.8860339012559212,0.8860339012559212,0.8860339012559212,1.211707852076548,0.9382214904797996,1.032096159763384,1.211707852076548,0.8860339012559212,0.9382214904797996,1.074964671278174,0.9895971763399439,1.011987747937423,1.063634068252145,1.011987747937423,1.019234026503089,0.9925940014536556,1.01479408181879)
$newRow5 = @(1.303066193037559,0.8523169836467658,0.8421598848881794,1.091285160017975,0.6139390195868498,0.6139390195868498,0.6139390195868498,1.467876293991917,0.9457490601414295,1.119862889078076,1.467876293991917,0.6139390195868498,0.9457490601414295,1.206812677066673,1.018517110079702,1.009188124573399,1.168303504717107,1.009188124573399,1.029712383434543,0.9465577106650043,1.029531935548594)

for ($i = 0; $i -lt $newRow4.Length; $i++) {
    $ws.Cells.Item(4, 3 + $i).Value = $newRow4[$i]
}
for ($i = 0; $i -lt $newRow5.Length; $i++) {
    $ws.Cells.Item(5, 3 + $i).Value = $newRow5[$i]
}

# ---------------------------------------------------------------------------
# 5) Rewrite columns A (method index) and B (method name) for every data row
#    (2..31) to match the new, longer method list. Re-apply the bordered
#    "index" style (same as A3) to the A cells that didn't already carry it
#    over from the shift above: the two brand-new rows (4:5) plus the two
#    rows that now exist past the sheet's old bottom edge (30:31), so no new
#    style slot gets allocated for them.
# ---------------------------------------------------------------------------
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A30:A31").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

for ($idx = 0; $idx -lt $methods.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 1).Value = $idx
    $ws.Cells.Item($row, 2).Value = $methods[$idx]
}

Write-Host "Inserted Holden/Rizzie Spiral, renamed Thomas Hex -> Matthies Hex, reran sweep."
